$d = $word.ActiveDocument
$wNs = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

$titleText = "Play Book of Tombs Online Slot for Free - Game Review"
$newMetaText = "Discover the Book of Tombs online slot game with 5 reels, 10 fixed paylines, high volatility rate, and free spins function. Play for free and read our review."

# -----------------------------------------------------------------------
# 1. Remove the "Meta description" paragraph that currently sits right
#    after the H1 title ("Play Book of Tombs Online Slot for Free -
#    Game Review").
# -----------------------------------------------------------------------
$metaIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -like "*Meta description*") {
        $metaIndex = $i
        break
    }
}
if ($metaIndex -ge 1) {
    $d.Paragraphs($metaIndex).Range.Delete()
}

# -----------------------------------------------------------------------
# 2. Find the paragraph that holds the "Create a cartoon-style..."
#    image-prompt text (now the final paragraph of the document) and
#    insert a new bold title paragraph right before it.
# -----------------------------------------------------------------------
$imageIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -like "*Create a cartoon-style feature image*") {
        $imageIndex = $i
        break
    }
}

$priorPara = $d.Paragraphs($imageIndex - 1)
$end = $priorPara.Range.Duplicate()
$end.Collapse(0)
$end.InsertParagraphAfter()

$newPara = $d.Paragraphs($imageIndex)
[void]$newPara.Range.InsertXML("<w:p $wNs><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>$titleText</w:t></w:r></w:p>")

# -----------------------------------------------------------------------
# 3. Replace the text of the (now shifted) image-prompt paragraph with
#    the meta description copy, keeping its existing italic run
#    formatting intact.
# -----------------------------------------------------------------------
$imageIndex2 = $imageIndex + 1
$targetPara = $d.Paragraphs($imageIndex2)
$targetRange = $targetPara.Range
$textOnly = $d.Range($targetRange.Start, $targetRange.End - 1)
$textOnly.Text = $newMetaText

Write-Output "done"
